$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H ("Absent") should be the complement of column E ("Real") for each
# attendance row (rows 3 through 21): a student is Absent when they were not
# marked as Real attendance for that date.
for ($r = 3; $r -le 21; $r++) {
    $real = $ws.Cells.Item($r, 5).Value()   # Column E = Real
    if ($real -eq 1) {
        $ws.Cells.Item($r, 8).Value = 0
    } else {
        $ws.Cells.Item($r, 8).Value = 1
    }
}
